# Apply "y" markers to previously blank cells in the phenology tracking grid.
# New cells that fall in an already-shaded column (style index 1) pick up
# that shading by copying the format from a neighboring cell in the same row.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GC_splist")

$xlPasteFormats = -4122

function Set-Y {
    param($cell)
    $ws.Range($cell).Value = "y"
}

function Set-Y-WithFormatFrom {
    param($cell, $formatSource)
    $ws.Range($formatSource).Copy()
    $ws.Range($cell).PasteSpecial($xlPasteFormats)
    $ws.Range($cell).Value = "y"
}

Set-Y "H2"
Set-Y "H3"

# Row 4: E4/F4/G4 already carry the shaded style (s=1); H4/I4 are brand
# new cells that need to inherit that same shading.
Set-Y "E4"
Set-Y "F4"
Set-Y "G4"
Set-Y-WithFormatFrom "H4" "G4"
Set-Y-WithFormatFrom "I4" "G4"

# Row 6: G6 keeps its shading (s=1) with no value; H6/I6 are new shaded cells.
Set-Y-WithFormatFrom "H6" "G6"
Set-Y-WithFormatFrom "I6" "G6"

Set-Y "H7"

# Row 8: G8 already shaded (s=1); H8 is new and needs the same shading.
Set-Y "G8"
Set-Y-WithFormatFrom "H8" "G8"

Set-Y "G9"
Set-Y "G10"

Set-Y "H11"
Set-Y "I11"

Set-Y "G12"

Set-Y "I13"

Set-Y "H14"

Set-Y "H15"

Set-Y "G16"
Set-Y "H16"

Set-Y "H18"

Set-Y "H21"

Set-Y "G22"

# Update the view: zoom to 140% and move the selection to H15
$ws.Activate()
$excel.ActiveWindow.Zoom = 140
$ws.Range("H15").Select()
